$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'ärgern',
    'spüren',
    'fahren',
    'planen',
    'räumen',
    'ändern',
    'kosten',
    'runden',
    'bellen',
    'tropfen',
    'platzen',
    'schwächen',
    'fällen',
    'arten',
    'bluten',
    'zünden',
    'dienen',
    'achten',
    'grüßen',
    'schlucken',
    'spinnen',
    'heilen',
    'erben',
    'irren',
    'decken',
    'sperren',
    'folgen',
    'flehen',
    'spielen',
    'boxen',
    'münzen',
    'betteln',
    'stehlen',
    'lockern',
    'quälen',
    'scheinen',
    'töten',
    'altern',
    'saufen',
    'kichern',
    'reizen',
    'drehen',
    'sterben',
    'hauen',
    'klingen',
    'klettern',
    'formen',
    'werden',
    'wachsen',
    'rufen',
    'tollen',
    'fischen',
    'enden',
    'pfeifen',
    'ehren',
    'führen',
    'machen',
    'wehtun',
    'trauen',
    'äußern',
    'scheitern',
    'schwören',
    'heulen',
    'biegen',
    'rasen',
    'liefern',
    'lügen',
    'freuen',
    'weichen',
    'seufzen',
    'bauen',
    'fangen',
    'jubeln',
    'werfen',
    'wenden',
    'brauchen',
    'lesen',
    'bergen',
    'wundern',
    'loben',
    'schrecken',
    'mauern',
    'flüchten',
    'gelten',
    'greifen',
    'pflanzen',
    'jagen',
    'wüten',
    'sorgen',
    'zögern',
    'helfen',
    'bitten',
    'zeigen',
    'kümmern',
    'fallen',
    'mögen',
    'dringen',
    'suchen',
    'kehren',
    'wirken',
    'stecken',
    'heben',
    'liegen',
    'graben',
    'schenken',
    'knarren',
    'sinken',
    'filmen',
    'zielen',
    'schwingen',
    'treiben',
    'schreiten',
    'siegen',
    'feiern',
    'geben',
    'hören',
    'gründen',
    'schmecken',
    'backen',
    'malen',
    'sichern',
    'schulden',
    'streichen',
    'sprengen',
    'warnen',
    'fließen',
    'trennen',
    'stammen'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $words[$i]
}
